$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("14:14").Copy()
$ws.Rows("15:15").Insert()
